$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Decoding Human Memory: Unraveling the Mysteries of the Mind" "Interlacing the Essence of Mathematics and Art: A Harmonic Tapestry"

# Author
Replace-Text "Alex Bishop" "Sophia Lancaster"

# Email (three separate runs)
Replace-Text "info" "SophiaLancaster05@edu"
Replace-Text "alexbishop@gmail" "sch"
Replace-Text "com" "uk"

# Paragraph 1 (intro) runs
Replace-Text "Throughout the annals of human history, the intricate workings of our memories have fascinated and perplexed philosophers, scientists, and artists alike" "Mathematics and art, seemingly dissimilar domains, converge in a symphony of harmony and interconnectedness"

Replace-Text " As we navigate the complexities of our lives, a tapestry of memories weaves its way through every aspect of our being, shaping our identities, guiding our decisions, and fueling our emotions" " Mathematics, with its abstract symbols and intricate equations, may appear distant from the vibrant strokes and expressive canvases of art"

Replace-Text " Yet, the enigmatic nature of memory continues to challenge our understanding" " However, beneath the surface, a profound kinship emerges, weaving together the fabric of intellectual exploration and creative expression"

Replace-Text " In this exploration, we embark on a journey to unravel the secrets of human memory, delving into the profound depths of the mind and investigating the diverse perspectives that shape our comprehension of this remarkable faculty" " This essay delves into the depths of this captivating relationship, revealing the profound synergy between mathematics and art"

Replace-Text "Unveiling the nature of memory unveils a realm of diverse facets and dimensions" "At the heart of these seemingly disparate disciplines lies a shared pursuit of patterns, structures, and relationships"

Replace-Text " From the fleeting impressions of short-term memory to the enduring archive of long-term storage, each facet offers a distinct glimpse into the intricacies of the mind's mechanisms" " Mathematics seeks to unravel the underlying order in the universe, while art endeavors to capture the essence of human emotion and experience"

Replace-Text " Moreover, the remarkable ability to recall and relive past experiences, the curious phenomenon of forgetting, and the enigmatic role of memory in our perception of time and reality paint a complex tapestry of interconnected processes" " In this convergence, both mathematics and art become lenses through which we interpret and express the enigmatic beauty of the world around us"

Replace-Text " As we traverse the terrain of memory, we discover the remarkable interplay between biology and psychology, as neural networks and cognitive structures intertwine to orchestrate this intricate cognitive landscape" " From the golden ratio's harmonious proportions to the intricate fractals that mimic the complexities of nature, mathematics and art intertwine, creating a kaleidoscope of patterns that fascinate and inspire"

Replace-Text "The quest to unravel the mysteries of human memory unfolds across disciplines, inviting a kaleidoscope of perspectives" "Moreover, both mathematics and art transcend the boundaries of language, appealing to a universal human experience"

Replace-Text " Neuroscience ventures into the realm of brain anatomy and neurochemistry, illuminating the physical basis of memory formation and retrieval" " Mathematical symbols and artistic forms speak a language of their own, capable of communicating concepts and emotions that words alone cannot capture"

Replace-Text " Psychology delves into the depths of cognitive processes, exploring the intricate mechanisms of encoding, storage, and recall" " They possess the unique ability to transcend cultural and linguistic barriers, resonating with individuals from all walks of life"

Replace-Text " In this convergence of disciplines, we forge a comprehensive understanding of memory, illuminating its profound impact on our lives and opening up new avenues for exploration and discovery" " In this shared capacity for transcendence, mathematics and art unite, becoming powerful tools for human connection and understanding"

# Summary heading + body
Replace-Text "The exploration of human memory unveils a realm of captivating complexity, where diverse facets and dimensions intertwine to create a multifaceted cognitive landscape" "Mathematics and art, though seemingly disparate disciplines, share a profound interconnectedness"

Replace-Text " Our understanding of this remarkable faculty emerges from the convergence of biological and psychological perspectives, revealing the intricate interplay between neural networks and cognitive structures" " Both seek to explore patterns, structures, and relationships, employing unique tools to unveil the enigmatic beauty of the world around us"

Replace-Text " As we delve deeper into this enigmatic realm, we uncover the mysteries of the mind, unlocking the secrets of memory and gaining a profound appreciation for this remarkable human capacity" " This essay has illuminated the rich tapestry woven by mathematics and art, highlighting their harmonious convergence as powerful tools for exploration, expression, and understanding"

Write-Output "done with simple replacements"
